$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Marche row (row 6): spatial checks completed, date moved forward to 12 Feb 2018 (serial 43143)
$ws.Range("B6").Value = 43143
$ws.Range("C6").Value = "Completato controlli qualità interni e controlli spaziali"

# Toscana row (row 7): same update
$ws.Range("B7").Value = 43143
$ws.Range("C7").Value = "Completato controlli qualità interni e controlli spaziali"

# Update the active cell selection to B8 (as recorded in the saved view state)
$ws.Range("B8").Select()
